# Forms the consolidated report: recompute the "Absent" column (H) for
# each attendance row as 1 when the "Real" attendance count (E) is 0,
# otherwise 0.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 3; $r -le $lastRow; $r++) {
    $realValue = $ws.Cells.Item($r, 5).Value2
    if ($realValue -eq 0) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
